$wb = $excel.ActiveWorkbook

# --- Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0" ---
$includeSheet = $wb.Worksheets.Item(2)
$includeSheet.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item(1)

# Update Date and Contact values in place (rows 8 and 10 keep their
# position; the new row is inserted further down at row 11).
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Shift rows 11-14 (Description, Purpose, Copyright, Immutable) down to
# 12-15 to make room for a new "Jurisdiction" row at 11. Copy bottom-up
# (with formatting) instead of using Rows.Insert so the existing cell
# style (s="2") is reused rather than new style indexes being created.
# Copy() does not blank out a destination cell when the source cell is
# empty, so ClearContents() the destination immediately beforehand.
$ws.Range("A15:B15").ClearContents()
$ws.Range("A14:B14").Copy($ws.Range("A15:B15"))

$ws.Range("A14:B14").ClearContents()
$ws.Range("A13:B13").Copy($ws.Range("A14:B14"))

$ws.Range("A13:B13").ClearContents()
$ws.Range("A12:B12").Copy($ws.Range("A13:B13"))

$ws.Range("A12:B12").ClearContents()
$ws.Range("A11:B11").Copy($ws.Range("A12:B12"))

# New Jurisdiction row (no published value).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
